$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'67.500.74"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -0.02%  "

$ws.Cells.Item(3, 4).Value = "'3.241.79"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -0.41%  "

$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.14%  "

$ws.Cells.Item(5, 4).Value = "'577.60"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.28%  "

$ws.Cells.Item(6, 4).Value = "'180.60"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.48%  "

$ws.Cells.Item(7, 4).Value = "'0.996"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.43%  "

$ws.Cells.Item(8, 4).Value = "'0.593"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -0.88%  "

$ws.Cells.Item(9, 4).Value = "'0.133"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +0.47%  "

$ws.Cells.Item(10, 4).Value = "'6.62"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -1.68%  "

$ws.Cells.Item(11, 4).Value = "'0.417"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.42%  "

$ws.Cells.Item(12, 4).Value = "'3.796.15"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -0.69%  "

$ws.Cells.Item(13, 5).Value = "  -0.34%  "

$ws.Cells.Item(14, 4).Value = "'28.04"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -1.26%  "

$ws.Cells.Item(15, 4).Value = "'67.891.95"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +0.57%  "

$ws.Cells.Item(16, 4).Value = "'0.0000170"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +1.18%  "

$ws.Cells.Item(17, 4).Value = "'3.226.15"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -1.10%  "

$ws.Cells.Item(18, 4).Value = "'5.79"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -0.77%  "

$ws.Cells.Item(19, 4).Value = "'13.43"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.58%  "

$ws.Cells.Item(20, 4).Value = "'390.96"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +4.27%  "

$ws.Cells.Item(21, 4).Value = "'7.65"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +0.36%  "

$ws.Cells.Item(22, 4).Value = "'1.00"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.03%  "

$ws.Cells.Item(23, 4).Value = "'70.61"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.35%  "

$ws.Cells.Item(24, 4).Value = "'0.515"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.87%  "

$ws.Cells.Item(25, 4).Value = "'0.0000119"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.66%  "

$ws.Cells.Item(26, 4).Value = "'0.187"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +3.37%  "

$ws.Cells.Item(27, 4).Value = "'9.50"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.86%  "

$ws.Cells.Item(28, 5).Value = "  -0.19%  "

$ws.Cells.Item(29, 5).Value = "  -0.75%  "

$ws.Cells.Item(30, 4).Value = "'5.62"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -0.79%  "

$ws.Cells.Item(31, 4).Value = "'22.91"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +0.89%  "

$ws.Cells.Item(32, 4).Value = "'7.08"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +2.43%  "

$ws.Cells.Item(33, 4).Value = "'0.998"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +0.03%  "

$ws.Cells.Item(34, 4).Value = "'1.27"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -0.02%  "

$ws.Cells.Item(35, 4).Value = "'163.56"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +0.13%  "

$ws.Cells.Item(36, 5).Value = "  -1.83%  "

$ws.Cells.Item(37, 4).Value = "'1.87"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +1.37%  "

$ws.Cells.Item(38, 2).Value = "Mantle"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(38, 4).Value = "'0.819"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -3.39%  "

$ws.Cells.Item(39, 2).Value = "EnergySwap"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(39, 4).Value = "'26.61"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -0.38%  "

$ws.Cells.Item(40, 4).Value = "'4.59"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -1.20%  "

$ws.Cells.Item(41, 4).Value = "'6.50"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -4.56%  "

$ws.Cells.Item(42, 2).Value = "dogwifhat"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(42, 4).Value = "'2.49"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -4.31%  "

$ws.Cells.Item(43, 2).Value = "OKB"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(43, 4).Value = "'41.30"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +1.29%  "

$ws.Cells.Item(44, 4).Value = "'0.0682"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.52%  "

$ws.Cells.Item(45, 4).Value = "'2.615.97"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -3.00%  "

$ws.Cells.Item(46, 4).Value = "'339.59"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -3.35%  "

$ws.Cells.Item(47, 4).Value = "'24.71"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -2.46%  "

$ws.Cells.Item(48, 5).Value = "  -0.54%  "

$ws.Cells.Item(49, 4).Value = "'6.30"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +2.57%  "

$ws.Cells.Item(50, 4).Value = "'0.102"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.38%  "

$ws.Cells.Item(51, 4).Value = "'31.21"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -0.34%  "
